# Scheduled-runner data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H, I, J, K, L, M, N) across the per-job Leve-profit
# sheets with freshly pulled market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 11999.5
$ws.Range("J3").Value = 11999.5
$ws.Range("L3").Value = 11999.5
$ws.Range("N3").Value = -12227.5
$ws.Range("H40").Value = 2400
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825
$ws.Range("H63").Value = 49900
$ws.Range("J63").Value = 49900
$ws.Range("L63").Value = 49900
$ws.Range("N63").Value = -51148
$ws.Range("H66").Value = 49900
$ws.Range("J66").Value = 49900
$ws.Range("L66").Value = 149700
$ws.Range("N66").Value = -155940
$ws.Range("H98").Value = 1785.1522
$ws.Range("I98").Value = 1446.5714
$ws.Range("J98").Value = 2862.4546
$ws.Range("K98").Value = 1446.5714
$ws.Range("L98").Value = 2862.4546
$ws.Range("M98").Value = 51.42859999999996
$ws.Range("N98").Value = -5858.4546
$ws.Range("H102").Value = 11999.5
$ws.Range("J102").Value = 11999.5
$ws.Range("L102").Value = 11999.5
$ws.Range("N102").Value = -18489.5
$ws.Range("H105").Value = 46599.6
$ws.Range("J105").Value = 46599.6
$ws.Range("L105").Value = 46599.6
$ws.Range("N105").Value = -53587.6
$ws.Range("H106").Value = 3799.8572
$ws.Range("I106").Value = 3376.3333
$ws.Range("J106").Value = 4562.2
$ws.Range("K106").Value = 3376.3333
$ws.Range("L106").Value = 4562.2
$ws.Range("M106").Value = -2745.3333
$ws.Range("N106").Value = -5824.2
$ws.Range("H122").Value = 1785.1522
$ws.Range("I122").Value = 1446.5714
$ws.Range("J122").Value = 2862.4546
$ws.Range("K122").Value = 4339.7142
$ws.Range("L122").Value = 8587.363799999999
$ws.Range("M122").Value = -1889.7142
$ws.Range("N122").Value = -13487.3638
$ws.Range("H132").Value = 1285.8695
$ws.Range("I132").Value = 1260.4286
$ws.Range("K132").Value = 3781.2858
$ws.Range("M132").Value = -1251.2858
$ws.Range("H138").Value = 2269.7012
$ws.Range("J138").Value = 2129.1177
$ws.Range("L138").Value = 6387.353099999999
$ws.Range("N138").Value = -16667.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3295.0815
$ws.Range("I32").Value = 2000.409
$ws.Range("J32").Value = 14688.2
$ws.Range("K32").Value = 2000.409
$ws.Range("L32").Value = 14688.2
$ws.Range("M32").Value = -1713.409
$ws.Range("N32").Value = -15262.2
$ws.Range("H97").Value = 1675.0454
$ws.Range("I97").Value = 1683.25
$ws.Range("K97").Value = 1683.25
$ws.Range("M97").Value = -1187.25
$ws.Range("H101").Value = 19998
$ws.Range("J101").Value = 19998
$ws.Range("L101").Value = 19998
$ws.Range("N101").Value = -26488
$ws.Range("H132").Value = 2427.3103
$ws.Range("I132").Value = 1872.75
$ws.Range("J132").Value = 5089.2
$ws.Range("K132").Value = 5618.25
$ws.Range("L132").Value = 15267.6
$ws.Range("M132").Value = -3088.25
$ws.Range("N132").Value = -20327.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
$ws.Range("H134").Value = 6769.2
$ws.Range("J134").Value = 2584
$ws.Range("L134").Value = 7752
$ws.Range("N134").Value = -12822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3447
$ws.Range("I122").Value = 1799.8
$ws.Range("K122").Value = 5399.4
$ws.Range("M122").Value = -2949.4
$ws.Range("H132").Value = 1837
$ws.Range("J132").Value = 3228
$ws.Range("L132").Value = 9684
$ws.Range("N132").Value = -14744
$ws.Range("H141").Value = 74148.25
$ws.Range("J141").Value = 74864.336
$ws.Range("L141").Value = 74864.336
$ws.Range("N141").Value = -85224.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 192.16667
$ws.Range("I2").Value = 229.6
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 1377.6
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = -1264.6
$ws.Range("N2").Value = -256
$ws.Range("H9").Value = 5249
$ws.Range("J9").Value = 8498
$ws.Range("L9").Value = 25494
$ws.Range("N9").Value = -25942
$ws.Range("H48").Value = 990
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H113").Value = 4292.241
$ws.Range("I113").Value = 25787.75
$ws.Range("K113").Value = 77363.25
$ws.Range("M113").Value = -75193.25
$ws.Range("H131").Value = 20583.17
$ws.Range("J131").Value = 21087.75
$ws.Range("L131").Value = 63263.25
$ws.Range("N131").Value = -73343.25
$ws.Range("H141").Value = 3363.5264
$ws.Range("I141").Value = 3053.353
$ws.Range("K141").Value = 9160.059000000001
$ws.Range("M141").Value = -3980.059000000001
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 58516.5
$ws.Range("J52").Value = 58516.5
$ws.Range("L52").Value = 58516.5
$ws.Range("N52").Value = -59034.5
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 6000
$ws.Range("N104").Value = -12988
$ws.Range("H113").Value = 1471.3572
$ws.Range("I113").Value = 1126.75
$ws.Range("J113").Value = 1930.8334
$ws.Range("K113").Value = 1126.75
$ws.Range("L113").Value = 1930.8334
$ws.Range("M113").Value = 1043.25
$ws.Range("N113").Value = -6270.8334
$ws.Range("H132").Value = 1604858.8
$ws.Range("I132").Value = 1833537.4
$ws.Range("J132").Value = 4109
$ws.Range("K132").Value = 5500612.199999999
$ws.Range("L132").Value = 12327
$ws.Range("M132").Value = -5498082.199999999
$ws.Range("N132").Value = -17387

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1980.3334
$ws.Range("I46").Value = 1554.1428
$ws.Range("K46").Value = 1554.1428
$ws.Range("M46").Value = -1366.1428
$ws.Range("H69").Value = 125000
$ws.Range("J69").Value = 125000
$ws.Range("L69").Value = 125000
$ws.Range("N69").Value = -126622
$ws.Range("H72").Value = 125000
$ws.Range("J72").Value = 125000
$ws.Range("L72").Value = 375000
$ws.Range("N72").Value = -383112
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 11200
$ws.Range("J103").Value = 11200
$ws.Range("L103").Value = 11200
$ws.Range("N103").Value = -13544
$ws.Range("H104").Value = 21249
$ws.Range("J104").Value = 21249
$ws.Range("L104").Value = 21249
$ws.Range("N104").Value = -28237
$ws.Range("H113").Value = 686
$ws.Range("I113").Value = 484.14285
$ws.Range("K113").Value = 1452.42855
$ws.Range("M113").Value = 717.5714499999999
$ws.Range("H136").Value = 26456694
$ws.Range("I136").Value = 39683400
$ws.Range("K136").Value = 119050200
$ws.Range("M136").Value = -119047650
